$d = $word.ActiveDocument

# Title (appears twice: document heading and bold "meta title" paragraph near the end)
$d.Content.Find.Execute("Play Bounty Showdown for Free - Review of Fantasma Games' Video Slot", $true, $false, $false, $false, $false, $true, 1, $false, "Play Bounty Showdown for Free", 2)

# "What we like" bullet list
$d.Content.Find.Execute("Exciting Double Bonus Respins and Showdown Free Spins features", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting Wild West theme", 2)
$d.Content.Find.Execute("Visually appealing Wild West-themed graphics", $true, $false, $false, $false, $false, $true, 1, $false, "Engaging bonus features", 2)
$d.Content.Find.Execute("High-quality sound and immersive gaming experience", $true, $false, $false, $false, $false, $true, 1, $false, "Visually appealing comic-style graphics", 2)
$d.Content.Find.Execute("Playable across desktop, tablet, and mobile devices", $true, $false, $false, $false, $false, $true, 1, $false, "High-quality sound effects", 2)

# "What we don't like" bullet list
$d.Content.Find.Execute("Limited number of paylines (10)", $true, $false, $false, $false, $false, $true, 1, $false, "Limited number of active paylines", 2)
$d.Content.Find.Execute("High volatility and potential for low payouts", $true, $false, $false, $false, $false, $true, 1, $false, "No progressive jackpot feature", 2)

# Meta description (italic paragraph)
$d.Content.Find.Execute("Read our review of Bounty Showdown by Fantasma Games and play this Wild West-themed online video slot for free. Learn about its features and payouts.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Bounty Showdown and play this exciting Wild West-themed slot game for free.", 2)
